$wb = $excel.ActiveWorkbook

# --- Sheet: P_valores ---
$wsP = $wb.Worksheets.Item("P_valores")

$wsP.Range("C2").Value = 0.1092544143289513
$wsP.Range("D2").Value = 0.1626591030465032
$wsP.Range("E2").Value = 0.884331627809273
$wsP.Range("F2").Value = 0.2404582752421955

$wsP.Range("B3").Value = 0.1092544143289513
$wsP.Range("D3").Value = 0.9617955425706066
$wsP.Range("E3").Value = 0.2182810110524067
$wsP.Range("F3").Value = 0.5771832901653418

$wsP.Range("B4").Value = 0.1626591030465032
$wsP.Range("C4").Value = 0.9617955425706066
$wsP.Range("E4").Value = 0.2710934703354031
$wsP.Range("F4").Value = 0.6555188217173904

$wsP.Range("B5").Value = 0.884331627809273
$wsP.Range("C5").Value = 0.2182810110524067
$wsP.Range("D5").Value = 0.2710934703354031
$wsP.Range("F5").Value = 0.08062213318340294

$wsP.Range("B6").Value = 0.2404582752421955
$wsP.Range("C6").Value = 0.5771832901653418
$wsP.Range("D6").Value = 0.6555188217173904
$wsP.Range("E6").Value = 0.08062213318340294

# --- Sheet: Estadisticos_DM ---
$wsDM = $wb.Worksheets.Item("Estadisticos_DM")

$wsDM.Range("C2").Value = -1.669167359367889
$wsDM.Range("D2").Value = -1.444614438951007
$wsDM.Range("E2").Value = -0.1471784609164602
$wsDM.Range("F2").Value = -1.206448970627254

$wsDM.Range("B3").Value = 1.669167359367889
$wsDM.Range("D3").Value = -0.04844902438631072
$wsDM.Range("E3").Value = 1.26735147725283
$wsDM.Range("F3").Value = 0.565908466678087

$wsDM.Range("B4").Value = 1.444614438951007
$wsDM.Range("C4").Value = 0.04844902438631072
$wsDM.Range("E4").Value = 1.128916032119154
$wsDM.Range("F4").Value = 0.4522469379577835

$wsDM.Range("B5").Value = 0.1471784609164602
$wsDM.Range("C5").Value = -1.26735147725283
$wsDM.Range("D5").Value = -1.128916032119154
$wsDM.Range("F5").Value = -1.831370871653119

$wsDM.Range("B6").Value = 1.206448970627254
$wsDM.Range("C6").Value = -0.565908466678087
$wsDM.Range("D6").Value = -0.4522469379577835
$wsDM.Range("E6").Value = 1.831370871653119
